$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D numeric-looking strings are preserved as text (not auto-converted to numbers)
$dCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D35', 'D36', 'D37', 'D38', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $dCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.419.81'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.805.73'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '1.008'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').Value = '306.89'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('D7').Value = '0.4521'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '0.3599'
$ws.Range('E8').Value = '  -1.58%  '
$ws.Range('D9').Value = '46.24'
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('D10').Value = '0.07070'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').Value = '0.8922'
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('D12').Value = '0.07814'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = '19.42'
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').Value = '1.806.35'
$ws.Range('E14').Value = '  -0.78%  '
$ws.Range('D15').Value = '5.290'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').Value = '6.323'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').Value = '85.29'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('D18').Value = '1.010'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').Value = '0.000008501'
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').Value = '1.007'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '26.467.84'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').Value = '14.18'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = '4.968'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '2.043.55'
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('D25').Value = '10.51'
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('D26').Value = '1.957'
$ws.Range('E26').Value = '  -1.26%  '
$ws.Range('D27').Value = '152.74'
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('D28').Value = '17.79'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').Value = '2.077'
$ws.Range('E29').Value = '  +3.57%  '
$ws.Range('D30').Value = '112.09'
$ws.Range('E30').Value = '  -0.95%  '
$ws.Range('D31').Value = '4.854'
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').Value = '0.08697'
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('D34').Value = '2.814'
$ws.Range('E34').Value = '  +10.92%  '
$ws.Range('D35').Value = '4.455'
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D36').Value = '0.7252'
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('D37').Value = '1.105'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').Value = '1.077'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('D40').Value = '2.912'
$ws.Range('E40').Value = '  +1.67%  '
$ws.Range('D41').Value = '0.05120'
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('D42').Value = '0.5098'
$ws.Range('E42').Value = '  +3.87%  '
$ws.Range('D43').Value = '6.777'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').Value = '0.1513'
$ws.Range('E44').Value = '  -3.63%  '
$ws.Range('D45').Value = '8.024'
$ws.Range('E45').Value = '  -1.46%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = '1.008'
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.4669'
$ws.Range('E47').Value = '  +1.62%  '
$ws.Range('D48').Value = '10.03'
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('D49').Value = '100.49'
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('D50').Value = '1.577'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('D51').Value = '0.05987'
$ws.Range('E51').Value = '  -0.13%  '
